# Updated cryptos list (Price + Volume(1h)) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.132.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.800.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5095"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3898"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07744"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.099"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.21%  "

$ws.Range("E11").Value = "  -2.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.320"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.003"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.798.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.275"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001072"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06581"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.19%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.170.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.240"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.427"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.007.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.52%  "

$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "

$ws.Range("E33").Value = "  -0.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.529"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07022"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.073"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02334"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.73%  "

$ws.Range("E38").Value = "  -0.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.006"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.49%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.155"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.301"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5903"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.722"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.200"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.904"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06743"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.93%  "
